$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row at row 39 (pushes the existing rows 39-88
# down to 40-89, including the previously-last row which becomes row 89).
$ws.Rows.Item(39).Insert()

# Populate the new row with the latest weekly price entry.
$ws.Range("A39").Value = 3
$ws.Range("B39").Value = "Femacal de La Calera"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 45117
$ws.Range("E39").Value = 5
$ws.Range("F39").Value = 100112022
$ws.Range("G39").Value = "Arveja Verde"
$ws.Range("H39").Value = "Perfection"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 24000
$ws.Range("M39").Value = 24000
$ws.Range("N39").Value = "`$/saco 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 960
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
